# Edit script: apply number format "#.#" (numFmtId 165) to column G data cells
# across all sheets, and update the more-precise point values for the
# GGEE_Masc and GGEE_Fem sheets (fixing the scraped rounding).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Elite_Masc", "Elite_Fem", "GGEE_Masc", "GGEE_Fem")
$lastRows = @(53, 22, 146, 48)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    $lastRow = $lastRows[$i]
    $ws.Range("G2:G$lastRow").NumberFormat = "#.#"
}

# Updated (higher-precision) Puntos values for GGEE_Masc (sheet3)
$ggeeMascVals = @(100,99.21899999999999,96.78100000000001,96.628,96.521,96.377,96.197,95.706,95.196,94.94199999999999,94.48999999999999,93.366,92.583,92.53,91.86799999999999,91.628,91.34,91.283,91.126,90.84399999999999,90.794,90.05500000000001,89.176,88.65900000000001,88.617,87.938,87.90300000000001,87.715,87.617,87.46599999999999,87.218,87.16200000000001,87.09,86.718,86.366,86.31399999999999,86.214,86.08499999999999,85.874,85.852,85.78400000000001,85.776,85.768,85.455,85.196,85.08199999999999,85.005,84.688,83.875,83.465,83.116,82.98,82.871,82.81399999999999,82.498,82.09,81.956,81.608,81.07299999999999,80.691,80.193,79.867,79.57299999999999,79.562,79.34699999999999,79.206,78.47,78.38500000000001,77.804,77.428,77.303,77.242,77.172,76.172,76.145,75.321,75.277,74.075,73.94,73.794,73.636,73.288,73.146,72.393,71.852,71.768,71.705,71.596,71.355,71.233,71.098,71.069,70.82299999999999,70.675,70.084,70.048,69.947,69.94499999999999,69.788,69.76300000000001,69.20699999999999,68.761,68.58,68.491,67.44499999999999,66.196,65.759,65.371,64.361,63.927,63.517,62.652,62.461,61.867,61.055,57.5,56.711,54.364,53.07,52.671,52.099,50.893,7.756,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$wsMasc = $wb.Worksheets.Item("GGEE_Masc")
for ($r = 0; $r -lt $ggeeMascVals.Length; $r++) {
    $wsMasc.Cells.Item($r + 2, 7).Value2 = $ggeeMascVals[$r]
}

# Updated (higher-precision) Puntos values for GGEE_Fem (sheet4)
$ggeeFemVals = @(100,98.559,95.154,90.98699999999999,89.746,89.28100000000001,85.84099999999999,85.021,84.889,84.617,81.88800000000001,80.553,77.917,75.422,74.398,73.25700000000001,71.13,69.51900000000001,69.35599999999999,67.014,66.967,66.239,65.99299999999999,65.93000000000001,65.07899999999999,63.313,63.173,62.001,61.401,60.36,58.179,57.511,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$wsFem = $wb.Worksheets.Item("GGEE_Fem")
for ($r = 0; $r -lt $ggeeFemVals.Length; $r++) {
    $wsFem.Cells.Item($r + 2, 7).Value2 = $ggeeFemVals[$r]
}

Write-Host "Edit complete"
